$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B8/C8 were blank; fill them in with the same date (25-Nov-2022, serial 44890)
$ws.Range("B8").Value = 44890
$ws.Range("C8").Value = 44890

# Give B8 the "date" variant of row 8's fill (fillId 9 + short-date number format),
# matching the style pattern already used for columns B/C on every other row.
$ws.Range("B8").NumberFormat = "mm-dd-yy"

# Copy that exact formatting onto C8 so both cells end up sharing a single style
# record instead of two separate (but identical) ones.
$ws.Range("B8").Copy()
$ws.Range("C8").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# The saved selection moved from D7 to C10
$excel.Goto($ws.Range("C10"))
